$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C","D","E","F","G","I","J","K","M","O")
$startRow = 2

$data = @(
    @(0.02707368271244093,0.1781656595681369,0.1446918898624894,1.341293567397457,0.002458821648212779,0.6863248066017675,0.1576000511679396,1.325493459160668,0.4201918276503207,3.214595581085661),
    @(0.02441208013344465,0.175510064981772,0.1448156765348578,1.352735958002803,0.002461705477286613,0.6923138821479746,0.159211077902885,1.174822085372625,0.3889817653281327,3.252798895395557),
    @(0.02276774562533745,0.1739327526466354,0.1449591373338279,1.360679540209958,0.002463569904438587,0.6965226120446388,0.1602924078220553,1.082087063608896,0.3698591012823584,3.278666616111337),
    @(0.02209515529386863,0.1733034440111467,0.1450345781379863,1.364147211643484,0.002464353316905006,0.6983711345455248,0.160756218523872,1.044243303852568,0.3620771160938787,3.289813217329197),
    @(0.02198332140446979,0.1731997628347273,0.1450481310452414,1.364736936941974,0.002464484831832046,0.6986861334039496,0.1608346323737067,1.037956216724069,0.3607855813774066,3.291700633796538),
    @(0.02275868496510469,0.1739242109786474,0.1449600859831968,1.360725373052851,0.00246358037403419,0.696547001872446,0.1602985691781598,1.081576902808393,0.3697541069154298,3.278814493463955),
    @(0.02615807490639099,0.177239004471474,0.1447205811407137,1.345048259912275,0.002459796578097596,0.6882794260010243,0.1581363992604849,1.273589462729603,0.4094224854082498,3.227267194752031),
    @(0.03274317009352501,0.1841589487111293,0.144785533660027,1.321599128498129,0.002453117170361375,0.6762925136985984,0.1546281691148614,1.648276724742686,0.4875156651878569,3.145352198015672),
    @(0.03753096969234093,0.189495862380852,0.1451585779840663,1.308832796755098,0.002448656782063541,0.6700751972394841,0.1524976190870078,1.922345035692445,0.5450586660216743,3.096912482171234),
    @(0.03969800359185172,0.1919780769068495,0.1453988165323707,1.303997249395053,0.002446723729229579,0.6678119494917567,0.1516256109774829,2.04674568498001,0.5712697399514184,3.077438285131109),
    @(0.04051700459581298,0.1929257907690101,0.1454999190275572,1.302306177525054,0.00244600546456598,0.6670363956145522,0.1513093958630023,2.093811569002469,0.5811997560011406,3.070433237536605),
    @(0.04034069007988705,0.1927213399726497,0.1454776945127421,1.302664147063659,0.002446159545531883,0.667199797018128,0.1513768756821534,2.083676993398228,0.5790609594126721,3.071925452857272),
    @(0.03976541569876701,0.1920558909317549,0.1454069313679867,1.303855315563112,0.002446664362177471,0.6677465096959665,0.1515993151674841,2.05061867829437,0.5720866020211162,3.076854565869382),
    @(0.03941283322468792,0.1916492916943469,0.1453649056262378,1.304603187056948,0.00244697536429789,0.6680920062682389,0.1517373889768869,2.030363970421661,0.5678151714784576,3.079921930361991),
    @(0.03738912567420982,0.1893347327357731,0.1451442960513951,1.309168392980659,0.002448785037913871,0.6702344961304405,0.1525565642551392,1.914209426133539,0.5433463610288811,3.09823680023112),
    @(0.03614481688494209,0.1879287132375111,0.1450270159124862,1.312218133643292,0.002449919756089691,0.6716937335626554,0.1530840083815619,1.842880319165147,0.5283440063553542,3.110129166988486),
    @(0.03542809363501931,0.1871251351687846,0.1449661997617611,1.314063728712078,0.002450581455624135,0.6725862185558,0.1533965270805027,1.801828022309508,0.5197183235679148,3.117210355318235),
    @(0.0351852475661758,0.1868539402755829,0.1449467493612175,1.314704316008729,0.002450807050169481,0.6728975234456556,0.1535039109123559,1.787924085340137,0.5167983935421034,3.119649284751858),
    @(0.0362773824839735,0.1880778563047016,0.145038813440852,1.311884015205962,0.002449798028159133,0.6715328906087876,0.1530269142532674,1.850476100348885,0.5299406973719982,3.108838252261819),
    @(0.03993443161414234,0.1922511396103204,0.1454274414223526,1.303501637463611,0.002446515712705719,0.6675837134624771,0.1515335993102305,2.060329864737639,0.5741350230530884,3.075396731509244),
    @(0.04231513919160079,0.1950237787911959,0.1457404675218292,1.298839689192647,0.002444450590142417,0.6654777868214765,0.1506392168756605,2.197235827559837,0.6030442832688721,3.055694466426502),
    @(0.04104537991180734,0.1935398612661885,0.1455680024944002,1.301253062737786,0.002445545480732947,0.6665582108953743,0.1511090940673263,2.124189824906864,0.5876126801985464,3.066012509528889),
    @(0.03621745379155072,0.1880104138884775,0.1450334591887454,1.312034782790278,0.002449853032203971,0.6716054409197341,0.153052697590855,1.847042185305725,0.5292188352864713,3.109421114499639),
    @(0.03097049632602022,0.1822423062222214,0.1447107892073518,1.327160315073385,0.002454845309154418,0.6790815117337061,0.1554988101070833,1.648276724742686,0.4875156651878569,3.165454264108519)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $col = $cols[$j]
        $ws.Range("$col$row").Value = $rowValues[$j]
    }
}
